$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1634
$ws.Range("J17").Value = 1634
$ws.Range("L17").Value = 4902
$ws.Range("N17").Value = -5238

$ws.Range("H55").Value = 1023
$ws.Range("I55").Value = 545.7273
$ws.Range("K55").Value = 545.7273
$ws.Range("M55").Value = -331.7273

$ws.Range("H76").Value = 3071.4285
$ws.Range("I76").Value = 2500
$ws.Range("K76").Value = 2500
$ws.Range("M76").Value = -2185

$ws.Range("H79").Value = 3071.4285
$ws.Range("I79").Value = 2500
$ws.Range("K79").Value = 2500
$ws.Range("M79").Value = -1408

$ws.Range("H111").Value = 1596.75
$ws.Range("J111").Value = 1499.6666
$ws.Range("L111").Value = 4498.9998
$ws.Range("N111").Value = -10632.9998

$ws.Range("H132").Value = 903.875
$ws.Range("I132").Value = 930.43335
$ws.Range("K132").Value = 2791.30005
$ws.Range("M132").Value = -261.3000499999998

$ws.Range("H138").Value = 4834.4688
$ws.Range("I138").Value = 3755.5186
$ws.Range("J138").Value = 5621.811
$ws.Range("K138").Value = 11266.5558
$ws.Range("L138").Value = 16865.433
$ws.Range("M138").Value = -6126.5558
$ws.Range("N138").Value = -27145.433

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5517.923
$ws.Range("I32").Value = 3291.4
$ws.Range("K32").Value = 3291.4
$ws.Range("M32").Value = -3004.4

$ws.Range("H61").Value = 1262
$ws.Range("I61").Value = 1283.8334
$ws.Range("K61").Value = 1283.8334
$ws.Range("M61").Value = -1071.8334

$ws.Range("H74").Value = 3083.1667
$ws.Range("I74").Value = 1199.8
$ws.Range("K74").Value = 1199.8
$ws.Range("M74").Value = -325.8

$ws.Range("H77").Value = 3083.1667
$ws.Range("I77").Value = 1199.8
$ws.Range("K77").Value = 5999
$ws.Range("M77").Value = -1631

$ws.Range("H110").Value = 9137.799999999999
$ws.Range("I110").Value = 9769.714
$ws.Range("K110").Value = 9769.714
$ws.Range("M110").Value = -7724.714

$ws.Range("H122").Value = 912430
$ws.Range("I122").Value = 1114025.5
$ws.Range("K122").Value = 3342076.5
$ws.Range("M122").Value = -3339626.5

$ws.Range("H132").Value = 815.5454999999999
$ws.Range("I132").Value = 718.2
$ws.Range("J132").Value = 1789
$ws.Range("K132").Value = 2154.6
$ws.Range("L132").Value = 5367
$ws.Range("M132").Value = 375.3999999999996
$ws.Range("N132").Value = -10427

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 1262
$ws.Range("I136").Value = 1283.8334
$ws.Range("K136").Value = 3851.5002
$ws.Range("M136").Value = -1301.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1788.3636
$ws.Range("I94").Value = 1519.1111
$ws.Range("K94").Value = 1519.1111
$ws.Range("M94").Value = -1068.1111

$ws.Range("H99").Value = 2885.7144
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2331.318
$ws.Range("I58").Value = 1292.5555
$ws.Range("K58").Value = 1292.5555
$ws.Range("M58").Value = -1089.5555

$ws.Range("H86").Value = 5057.143
$ws.Range("I86").Value = 3400
$ws.Range("K86").Value = 3400
$ws.Range("M86").Value = -2277

$ws.Range("H89").Value = 5057.143
$ws.Range("I89").Value = 3400
$ws.Range("K89").Value = 17000
$ws.Range("M89").Value = -11384

$ws.Range("H103").Value = 29749.75
$ws.Range("I103").Value = 29749.75
$ws.Range("K103").Value = 29749.75
$ws.Range("M103").Value = -28577.75

$ws.Range("H122").Value = 4983.5
$ws.Range("I122").Value = 6225.5
$ws.Range("K122").Value = 18676.5
$ws.Range("M122").Value = -16226.5

$ws.Range("H134").Value = 2875.35
$ws.Range("I134").Value = 2783.1875
$ws.Range("K134").Value = 8349.5625
$ws.Range("M134").Value = -5814.5625

$ws.Range("H136").Value = 2331.318
$ws.Range("I136").Value = 1292.5555
$ws.Range("K136").Value = 3877.6665
$ws.Range("M136").Value = -1327.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3172334.8
$ws.Range("I4").Value = 5342352.5
$ws.Range("K4").Value = 16027057.5
$ws.Range("M4").Value = -16026945.5

$ws.Range("H106").Value = 12499.667
$ws.Range("J106").Value = 12499.667
$ws.Range("L106").Value = 37499.001
$ws.Range("N106").Value = -39391.001

$ws.Range("H121").Value = 783.6923
$ws.Range("I121").Value = 766.3333
$ws.Range("J121").Value = 788.9
$ws.Range("K121").Value = 2298.9999
$ws.Range("L121").Value = 2366.7
$ws.Range("M121").Value = -988.9998999999998
$ws.Range("N121").Value = -4986.7

$ws.Range("H131").Value = 1338.585
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1338.585
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4015.755
$ws.Range("N131").Value = -14095.755
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 60
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H80").Value = 11596.308
$ws.Range("I80").Value = 5234
$ws.Range("J80").Value = 17049.715
$ws.Range("K80").Value = 5234
$ws.Range("L80").Value = 17049.715
$ws.Range("M80").Value = -4236
$ws.Range("N80").Value = -19045.715

$ws.Range("H83").Value = 11596.308
$ws.Range("I83").Value = 5234
$ws.Range("J83").Value = 17049.715
$ws.Range("K83").Value = 26170
$ws.Range("L83").Value = 85248.575
$ws.Range("M83").Value = -21178
$ws.Range("N83").Value = -95232.575

$ws.Range("H102").Value = 1867.1428
$ws.Range("I102").Value = 1695
$ws.Range("J102").Value = 2900
$ws.Range("K102").Value = 1695
$ws.Range("L102").Value = 2900
$ws.Range("M102").Value = -73
$ws.Range("N102").Value = -6144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4198.7
$ws.Range("I40").Value = 3799.6
$ws.Range("J40").Value = 4597.8
$ws.Range("K40").Value = 3799.6
$ws.Range("L40").Value = 4597.8
$ws.Range("M40").Value = -3663.6
$ws.Range("N40").Value = -4869.8

$ws.Range("H43").Value = 1463522
$ws.Range("I43").Value = 500000
$ws.Range("J43").Value = 2137987.5
$ws.Range("K43").Value = 500000
$ws.Range("L43").Value = 2137987.5
$ws.Range("M43").Value = -499807
$ws.Range("N43").Value = -2138373.5

$ws.Range("H46").Value = 4701.231
$ws.Range("I46").Value = 2967
$ws.Range("K46").Value = 2967
$ws.Range("M46").Value = -2779

$ws.Range("H132").Value = 3385.375
$ws.Range("I132").Value = 3383.4285
$ws.Range("K132").Value = 10150.2855
$ws.Range("M132").Value = -7620.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6500.5
$ws.Range("I126").Value = 4751
$ws.Range("K126").Value = 14253
$ws.Range("M126").Value = -11783

$ws.Range("H132").Value = 2761.389
$ws.Range("I132").Value = 682.6
$ws.Range("K132").Value = 2047.8
$ws.Range("M132").Value = 482.1999999999998
